# "added colors to rows" - highlight the weekday DTR rows (5-8 and 11-15)
# with a red fill, bump the SICK LEAVE column to 1 for those rows, flip the
# stray " " placeholder in B19 to a real boolean FALSE, and drop the extra
# (unused) third argument from the FLOOR() calls in the summary formulas.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Color the two weekday blocks -------------------------------------------------
# (Sat/Sun rows 9-10 are intentionally left alone.)
$rngColor = $ws.Range("A5:J8,A11:J15")

# Keep the existing center/center/wrap-text look while (re)asserting it so the
# style carries the alignment flags explicitly.
$rngColor.HorizontalAlignment = -4108
$rngColor.VerticalAlignment = -4108
$rngColor.WrapText = $true

# Cycle through the palette that was tried before red was settled on.
$rngColor.Interior.Color = 13411113   # FF29A3CC - cyan (tried first)
$rngColor.Interior.Color = 6737151    # FFFFCC66 - orange (tried second)
$rngColor.Interior.Color = 6184671    # FFDF5E5E - red (final)

# --- Mark those days as a sick-leave day (column I) -------------------------------
foreach ($r in 5,6,7,8,11,12,13,14,15) {
    $ws.Cells.Item($r, 9).Value = 1
}

# --- B19 becomes a real boolean FALSE instead of a blank-space placeholder --------
$ws.Range("B19").Value = $false

# --- Simplify the FLOOR() calls: drop the redundant significance/mode arg --------
$ws.Range("B22").Formula = '=FLOOR(F17/8,1)&"."&FLOOR(MOD(F17,8),1)&"."&(MOD(F17,8)-FLOOR(MOD(F17,8),1))*60'
$ws.Range("B23").Formula = '=FLOOR(H19,1)&"."&(H19-FLOOR(H19,1))*8&".0"'
$ws.Range("B24").Formula = '=FLOOR(I19,1)&"."&(I19-FLOOR(I19,1))*8&".0"'
$ws.Range("B27").Formula = '=FLOOR(K27/8,1)&"."&FLOOR(MOD(K27,8),1)&"."&(MOD(K27,8)-FLOOR(MOD(K27,8),1))*60'
